# Daily attendance processing - 2025-10-29 23:20:46
# For every row in column G ("Recorded By"), when the cell contains a
# comma-separated list of at least two recorders and does not involve
# "admin@admin.com", swap the first two entries in the list (e.g.
# "dnasr281@gmail.com, System" <-> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    if ($val -like "*admin@admin.com*") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Length -ge 2) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value = $parts -join ", "
    }
}
